# Atualizado por script em 20-12-2023 02:45
#
# This script:
#  1. Rotates the F:V data of rows 58-60 (the match data was re-scraped and
#     re-ordered): new row58 <- old row60, new row59 <- old row58, new row60 <- old row59
#  2. Appends three new match rows (64, 65, 66) with formatting copied from
#     an existing data row so the new A/E columns keep their original styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: capture current (pre-edit) F:V values of rows 58, 59, 60 before
# overwriting anything, then write them back in rotated order.
# ---------------------------------------------------------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$old58 = @{}
$old59 = @{}
$old60 = @{}
foreach ($col in $cols) {
    $old58[$col] = $ws.Range("$col" + "58").Value()
    $old59[$col] = $ws.Range("$col" + "59").Value()
    $old60[$col] = $ws.Range("$col" + "60").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col" + "58").Value = $old60[$col]
    $ws.Range("$col" + "59").Value = $old58[$col]
    $ws.Range("$col" + "60").Value = $old59[$col]
}

# ---------------------------------------------------------------------------
# Step 2: append three new rows (64-66) at the bottom of the table.
# Copy formats from row 58 (A58 bold/boxed style, E58 date style) so the new
# rows match the sheet's existing look.
# ---------------------------------------------------------------------------
$ws.Range("A58:V58").Copy()
$ws.Range("A64:V66").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# NOTE: this runtime's PowerShell function-parameter binding does not
# reliably bind *named* parameters (-Foo bar), so positional parameters are
# used instead: Row, Indice, Pais, Torneio, Temporada, DataPartida,
# Home, HomeGols, Away, AwayGols,
# HomeOpen, HomeOpenDT, HomeClose, HomeCloseDT,
# DrawOpen, DrawOpenDT, DrawClose, DrawCloseDT,
# AwayOpen, AwayOpenDT, AwayClose, AwayCloseDT, Url
function Set-MatchRow {
    param(
        [int]$Row,
        [int]$Indice,
        [string]$Pais,
        [string]$Torneio,
        [string]$Temporada,
        [double]$DataPartida,
        [string]$Home,
        [int]$HomeGols,
        [string]$Away,
        [int]$AwayGols,
        [double]$HomeOpen,
        [string]$HomeOpenDT,
        [double]$HomeClose,
        [string]$HomeCloseDT,
        [double]$DrawOpen,
        [string]$DrawOpenDT,
        [double]$DrawClose,
        [string]$DrawCloseDT,
        [double]$AwayOpen,
        [string]$AwayOpenDT,
        [double]$AwayClose,
        [string]$AwayCloseDT,
        [string]$Url
    )

    $ws.Range("A$Row").Value = $Indice
    $ws.Range("B$Row").Value = $Pais
    $ws.Range("C$Row").Value = $Torneio
    $ws.Range("D$Row").Value = $Temporada
    $ws.Range("E$Row").Value = $DataPartida
    $ws.Range("F$Row").Value = $Home
    $ws.Range("G$Row").Value = $HomeGols
    $ws.Range("H$Row").Value = $Away
    $ws.Range("I$Row").Value = $AwayGols
    $ws.Range("J$Row").Value = $HomeOpen
    $ws.Range("K$Row").Value = $HomeOpenDT
    $ws.Range("L$Row").Value = $HomeClose
    $ws.Range("M$Row").Value = $HomeCloseDT
    $ws.Range("N$Row").Value = $DrawOpen
    $ws.Range("O$Row").Value = $DrawOpenDT
    $ws.Range("P$Row").Value = $DrawClose
    $ws.Range("Q$Row").Value = $DrawCloseDT
    $ws.Range("R$Row").Value = $AwayOpen
    $ws.Range("S$Row").Value = $AwayOpenDT
    $ws.Range("T$Row").Value = $AwayClose
    $ws.Range("U$Row").Value = $AwayCloseDT
    $ws.Range("V$Row").Value = $Url
}

Set-MatchRow 64 63 "israel" "ligat-ha-al" "2023-2024" `
    45279.75 `
    "Hapoel Petah Tikva" 0 "Maccabi Bnei Raina" 1 `
    2.69 "18/12/2023 12:42" 3 "19/12/2023 17:55" `
    3.04 "18/12/2023 12:42" 2.99 "19/12/2023 16:00" `
    2.69 "18/12/2023 12:42" 2.66 "19/12/2023 17:55" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-petah-tikva-maccabi-bnei-raina/jXHOdLVP/"

Set-MatchRow 65 64 "israel" "ligat-ha-al" "2023-2024" `
    45279.77083333334 `
    "Hapoel Jerusalem" 1 "Hapoel Hadera" 0 `
    1.99 "16/12/2023 18:13" 2.55 "19/12/2023 18:27" `
    3.26 "16/12/2023 18:13" 3.06 "19/12/2023 18:27" `
    3.79 "16/12/2023 18:13" 3.07 "19/12/2023 18:27" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-jerusalem-hapoel-hadera/rP9jhuvs/"

Set-MatchRow 66 65 "israel" "ligat-ha-al" "2023-2024" `
    45279.79166666666 `
    "Hapoel Tel Aviv" 2 "Maccabi Petah Tikva" 0 `
    2.37 "18/12/2023 12:42" 2.24 "19/12/2023 18:58" `
    3.45 "18/12/2023 12:42" 3.23 "19/12/2023 18:58" `
    2.77 "18/12/2023 12:42" 3.47 "19/12/2023 18:58" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-tel-aviv-maccabi-petah-tikva/CtJKc1GJ/"

Write-Output ("Dimension now: " + $ws.UsedRange.Address())
